# Updated cryptos list - price and volume(1h) refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.064.71"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").Value = "'1.716.54"
$ws.Range("E3").Value = "  +1.28%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.74%  "
$ws.Range("D5").Value = "'318.49"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").Value = "'1.007"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("D7").Value = "'0.3973"
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").Value = "'0.4133"
$ws.Range("E8").Value = "  +1.35%  "
$ws.Range("D9").Value = "'1.532"
$ws.Range("E9").Value = "  +2.51%  "
$ws.Range("D10").Value = "'1.009"
$ws.Range("E10").Value = "  +0.82%  "
$ws.Range("D11").Value = "'52.63"
$ws.Range("E11").Value = "  +2.35%  "
$ws.Range("D12").Value = "'0.08936"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").Value = "'7.721"
$ws.Range("E13").Value = "  +7.34%  "
$ws.Range("D14").Value = "'25.01"
$ws.Range("E14").Value = "  +6.57%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.00001395"
$ws.Range("E15").Value = "  +4.53%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'8.156"
$ws.Range("E16").Value = "  -0.87%  "

$ws.Range("D17").Value = "'1.710.34"
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("D18").Value = "'100.87"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D19").Value = "'0.07147"
$ws.Range("E19").Value = "  +2.03%  "
$ws.Range("E20").Value = "  +3.00%  "
$ws.Range("D21").Value = "'7.486"
$ws.Range("E21").Value = "  +6.29%  "
$ws.Range("D22").Value = "'1.009"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").Value = "'14.52"
$ws.Range("E23").Value = "  +2.10%  "
$ws.Range("D24").Value = "'25.068.17"
$ws.Range("E24").Value = "  +1.59%  "
$ws.Range("D25").Value = "'3.109"
$ws.Range("E25").Value = "  -0.68%  "
$ws.Range("D26").Value = "'2.355"
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("D27").Value = "'23.09"
$ws.Range("E27").Value = "  +1.81%  "
$ws.Range("D28").Value = "'165.67"
$ws.Range("E28").Value = "  +1.55%  "
$ws.Range("D29").Value = "'8.814"
$ws.Range("E29").Value = "  +19.38%  "
$ws.Range("D30").Value = "'139.49"
$ws.Range("E30").Value = "  +1.44%  "
$ws.Range("D31").Value = "'5.222"
$ws.Range("E31").Value = "  +1.03%  "
$ws.Range("D32").Value = "'7.813"
$ws.Range("E32").Value = "  +9.55%  "
$ws.Range("D33").Value = "'1.897.48"
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("D34").Value = "'0.09008"
$ws.Range("E34").Value = "  +4.68%  "
$ws.Range("D35").Value = "'1.080"
$ws.Range("E35").Value = "  +0.62%  "
$ws.Range("D36").Value = "'0.02997"
$ws.Range("E36").Value = "  +9.88%  "
$ws.Range("D37").Value = "'0.2794"
$ws.Range("E37").Value = "  +1.71%  "
$ws.Range("D38").Value = "'1.977"
$ws.Range("E38").Value = "  +2.77%  "
$ws.Range("D39").Value = "'11.06"
$ws.Range("E39").Value = "  -3.87%  "
$ws.Range("D40").Value = "'14.66"
$ws.Range("E40").Value = "  +1.57%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.8200"
$ws.Range("E41").Value = "  +6.74%  "

$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "'0.09263"
$ws.Range("E42").Value = "  +0.35%  "

$ws.Range("D43").Value = "'1.484"
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("D44").Value = "'16.53"
$ws.Range("E44").Value = "  +3.76%  "
$ws.Range("D45").Value = "'0.7410"
$ws.Range("E45").Value = "  +3.22%  "
$ws.Range("D46").Value = "'2.644"
$ws.Range("E46").Value = "  +0.81%  "
$ws.Range("D47").Value = "'4.291"
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("D48").Value = "'1.007"
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("D49").Value = "'1.351"
$ws.Range("E49").Value = "  +1.80%  "
$ws.Range("D50").Value = "'140.68"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("E51").Value = "  +4.49%  "
